# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet, a new (blank) column is inserted
# immediately before column N ("Late"), pushing the existing "Late",
# "heading"/Outstanding-original and Outstanding columns one slot to the
# right (N->O, O->P, P->Q). The new column picks up the same width as the
# column immediately to its left ("Principal Paid", column M).
#
# The active sheet also moves from "Transactions" to "Repayment schedule",
# with a new selected cell on each sheet.

$wb = $excel.ActiveWorkbook
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsTx = $wb.Worksheets.Item("Transactions")

# Insert a blank column before column N (14), matching the width of the
# column to its left (M).
$leftWidth = $wsRepay.Columns.Item(13).ColumnWidth
$wsRepay.Columns.Item(14).Insert()
$wsRepay.Columns.Item(14).ColumnWidth = $leftWidth

# Leave behind a specific selection on the Transactions sheet ...
$wsTx.Range("D10").Select()

# ... then make "Repayment schedule" the active sheet/tab with its own
# selection.
$wsRepay.Activate()
$wsRepay.Range("J13").Select()
